# Apply the commit: "add docs, apk build"
#  1. Update the datetimeFigureOut footer field cached text on the
#     slide master and every slide layout (2018-11-30 -> 2018-12-06).
#  2. Append a brand-new slide 3 containing a "drinking log" card:
#     a rounded card frame, a header bar, several labelled fields,
#     and four thin divider lines layered on top.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Refresh the footer date placeholder everywhere it appears.
# ---------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "2018-12-06"
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes "2018-12-06"
}

# ---------------------------------------------------------------
# 2. Build the new third slide.
# ---------------------------------------------------------------
$s = $p.Slides.Add(3, 12)

# -- card background frame --------------------------------------
$rect3 = $s.Shapes.AddShape(1, 4246789, 473527, 3698421, 5690507)
$rect3.Name = "직사각형 3"
$rect3.Fill.ForeColor.ObjectThemeColor = 2
$rect3.Line.ForeColor.ObjectThemeColor = 10

# -- header bar ----------------------------------------------------
$rect4 = $s.Shapes.AddShape(1, 4246789, 473527, 3698421, 579665)
$rect4.Name = "직사각형 4"
$rect4.Fill.ForeColor.ObjectThemeColor = 5
$rect4.Line.ForeColor.ObjectThemeColor = 5

# -- title bar: 술이름 + 먹은날짜 ----------------------------------
$rect5 = $s.Shapes.AddShape(1, 5152366, 551088, 2225379, 506358)
$rect5.Name = "직사각형 5"
$rect5.Fill.ForeColor.ObjectThemeColor = 6
$rect5.Line.ForeColor.ObjectThemeColor = 5
$rect5.TextFrame.TextRange.Text = "술이름 + 먹은날짜"
$rect5.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# -- 술이름: @@@ bar ------------------------------------------------
$rect6 = $s.Shapes.AddShape(1, 4352923, 1130725, 3539557, 607942)
$rect6.Name = "직사각형 6"
$rect6.Fill.ForeColor.ObjectThemeColor = 6
$rect6.Line.ForeColor.ObjectThemeColor = 5
$rect6.TextFrame.TextRange.Text = "술이름: @@@"
$rect6.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# -- 날짜 label ------------------------------------------------------
$rect7 = $s.Shapes.AddShape(1, 7359967, 2118585, 1080535, 245289)
$rect7.Name = "직사각형 7"
$rect7.Fill.ForeColor.ObjectThemeColor = 6
$rect7.Line.ForeColor.ObjectThemeColor = 5
$rect7.TextFrame.TextRange.Text = "날짜"
$rect7.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# -- 술종류 label ------------------------------------------------------
$rect8 = $s.Shapes.AddShape(1, 7359967, 2463049, 1080535, 245289)
$rect8.Name = "직사각형 8"
$rect8.Fill.ForeColor.ObjectThemeColor = 6
$rect8.Line.ForeColor.ObjectThemeColor = 5
$rect8.TextFrame.TextRange.Text = "술종류"
$rect8.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# -- big review block: 맛있는 술이였나요? --------------------------
$rect9 = $s.Shapes.AddShape(1, 4352924, 2807417, 3539557, 4458652)
$rect9.Name = "직사각형 9"
$rect9.Fill.ForeColor.ObjectThemeColor = 6
$rect9.Line.ForeColor.ObjectThemeColor = 5
$rect9.TextFrame.TextRange.Text = "맛있는 술이였나요? @@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@"
$rect9.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# -- 같이먹은 안주: @@@@ bar -----------------------------------------
$rect10 = $s.Shapes.AddShape(1, 4352923, 2118585, 2380316, 606246)
$rect10.Name = "직사각형 10"
$rect10.Fill.ForeColor.ObjectThemeColor = 6
$rect10.Line.ForeColor.ObjectThemeColor = 5
$rect10.TextFrame.TextRange.Text = "같이먹은 안주: @@@@"
$rect10.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# -- small corner square, top right --------------------------------
$rect11 = $s.Shapes.AddShape(1, 8280077, 551088, 422563, 506358)
$rect11.Name = "직사각형 11"
$rect11.Fill.ForeColor.ObjectThemeColor = 4
$rect11.Line.ForeColor.ObjectThemeColor = 5

# -- small corner square, top left -----------------------------------
$rect12 = $s.Shapes.AddShape(1, 4352923, 551088, 424917, 506358)
$rect12.Name = "직사각형 12"
$rect12.Fill.ForeColor.ObjectThemeColor = 4
$rect12.Line.ForeColor.ObjectThemeColor = 5

# -- group all ten rectangles together into the card --------------
$cardRange = $s.Shapes.Range(@($rect3.Name, $rect4.Name, $rect5.Name, $rect6.Name, $rect7.Name, $rect8.Name, $rect9.Name, $rect10.Name, $rect11.Name, $rect12.Name))
$group = $cardRange.Group()
$group.Name = "그룹 13"

# -- divider lines sitting on top of the group ----------------------
$line15 = $s.Shapes.AddLine(5039239, 1693770, 9668389, 1693770)
$line15.Name = "직선 연결선 15"
$line15.Line.ForeColor.ObjectThemeColor = 5

$line16 = $s.Shapes.AddLine(5039239, 2281600, 9668389, 2281600)
$line16.Name = "직선 연결선 16"
$line16.Line.ForeColor.ObjectThemeColor = 5

$line18 = $s.Shapes.AddLine(7997980, 1565184, 7997980, 2422434)
$line18.Name = "직선 연결선 18"
$line18.Line.ForeColor.ObjectThemeColor = 5

$line20 = $s.Shapes.AddLine(7958404, 1985645, 9341532, 1985645)
$line20.Name = "직선 연결선 20"
$line20.Line.ForeColor.ObjectThemeColor = 5

Write-Host "Slide count now: " $p.Slides.Count
